$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow($Sheet, $Row, $D, $I, $J, $K, $L, $M, $P) {
    $Sheet.Range("A$Row").Value = 9
    $Sheet.Range("B$Row").Value = "Vega Central Mapocho de Santiago"
    $Sheet.Range("C$Row").Value = "Metropolitana"
    $Sheet.Range("D$Row").Value = $D
    $Sheet.Range("D$Row").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $Sheet.Range("E$Row").Value = 13
    $Sheet.Range("F$Row").Value = 100112035
    $Sheet.Range("G$Row").Value = "Bruselas (repollito)"
    $Sheet.Range("H$Row").Value = "Sin especificar"
    $Sheet.Range("I$Row").Value = $I
    $Sheet.Range("J$Row").Value = $J
    $Sheet.Range("K$Row").Value = $K
    $Sheet.Range("L$Row").Value = $L
    $Sheet.Range("M$Row").Value = $M
    $Sheet.Range("N$Row").Value = "`$/malla 15 kilos"
    $Sheet.Range("O$Row").Value = "Hijuelas"
    $Sheet.Range("P$Row").Value = $P
    $Sheet.Range("Q$Row").Value = 15
    $Sheet.Range("R$Row").Value = "Hortaliza"
}

# New weekly price observation inserted at the top of the data block (row 5).
# Existing rows 5-16 shift down to 6-17.
$ws.Rows.Item(5).Insert()
Set-DataRow $ws 5 44425 "Primera" 25 24000 25000 24520 1635

# Another new observation inserted after the row now holding the 44413 entry,
# i.e. before the row now holding the 44400 entry (originally row 11, now row 12).
$ws.Rows.Item(12).Insert()
Set-DataRow $ws 12 44421 "Primera" 18 24000 25000 24500 1633

# Final new observation appended as a new row at the end of the table.
Set-DataRow $ws 19 44418 "Primera" 16 25000 26000 25500 1700
